$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row with the additional test case value
$ws.Range("A4").Value = "Documents\no_classification.pdf"

# Update the selection to match the edited workbook's saved cursor position
$ws.Range("H13").Select()
